# edit.ps1 -- apply the functional-specs edit described by the diff:
#   1) "...offers filter result page" -> split into
#        "...offers " | "search" | [[_GoBack bookmark]] | " result page"
#   2) merge the four "manage functionality..." runs (which also had the
#      _GoBack bookmark) into a single run, dropping the bookmark there.

$d = $word.ActiveDocument

# Small helper-less pattern: to force the engine to recompute xml:space
# ("preserve") for a range based on its *actual* final text (instead of
# inheriting it from a larger, previously-merged run), we always write a
# short throw-away placeholder first and then the real text as a second
# assignment.

# =====================================================================
# Change 1: "One big map should be displayed on the offers filter result
#            page" -> "...offers search result page" with the run split
#            into three pieces and the _GoBack bookmark sitting right
#            after "search".
# =====================================================================

$find1 = $d.Content
$find1.Find.ClearFormatting()
$found1 = $find1.Find.Execute("One big map should be displayed on the offers filter result page", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $paraStart = $find1.Start

    # Step 1: swap "filter" -> "search" across the whole sentence (two-step
    # write so the resulting merged run's xml:space reflects its real text).
    $find1.Text = "PLACEHOLDERTEXT"
    $finalText = "One big map should be displayed on the offers search result page"
    $find1.Text = $finalText

    $offersEndOffset = $paraStart + $finalText.IndexOf("search")
    $searchEndOffset = $offersEndOffset + "search".Length
    $pageEndOffset   = $find1.End
    $dotEndOffset    = $pageEndOffset + 1

    # Step 2: drop temporary "barrier" bookmarks at every boundary we want
    # to keep as a distinct run, plus the real _GoBack bookmark right after
    # "search". Bookmarks stop adjacent same-format runs from being
    # re-coalesced once we touch the paragraph again.
    $d.Bookmarks.Add("ZZBarrierOffers", $d.Range($offersEndOffset, $offersEndOffset)) | Out-Null
    $d.Bookmarks.Add("ZZBarrierPage",   $d.Range($pageEndOffset, $pageEndOffset)) | Out-Null
    $d.Bookmarks.Add("ZZBarrierDot",    $d.Range($dotEndOffset, $dotEndOffset)) | Out-Null
    $d.Bookmarks.Add("_GoBack",         $d.Range($searchEndOffset, $searchEndOffset)) | Out-Null

    # Step 3: re-assert each sub-run's own text (two-step) so xml:space is
    # computed from that run's real content, right to left so offsets to
    # the left stay valid while we edit.
    $rDot = $d.Range($pageEndOffset, $dotEndOffset)
    $rDot.Text = "q1q1q1q1"
    $d.Range($pageEndOffset, $pageEndOffset + 8).Text = "."

    $rTail = $d.Range($searchEndOffset, $pageEndOffset)
    $rTail.Text = "q2q2q2q2"
    $d.Range($searchEndOffset, $searchEndOffset + 8).Text = " result page"

    $rSearch = $d.Range($offersEndOffset, $searchEndOffset)
    $rSearch.Text = "q3q3q3q3"
    $d.Range($offersEndOffset, $offersEndOffset + 8).Text = "search"

    $rHead = $d.Range($paraStart, $offersEndOffset)
    $rHead.Text = "q4q4q4q4"
    $d.Range($paraStart, $paraStart + 8).Text = "One big map should be displayed on the offers "

    # Step 4: drop the temporary barrier bookmarks, keep _GoBack.
    $d.Bookmarks("ZZBarrierOffers").Delete()
    $d.Bookmarks("ZZBarrierPage").Delete()
    $d.Bookmarks("ZZBarrierDot").Delete()
}

# =====================================================================
# Change 2: merge the "manage functionality" sentence's four runs (three
#           text runs plus the _GoBack bookmark) into a single run.
# =====================================================================

$find2 = $d.Content
$find2.Find.ClearFormatting()
$quote1 = [char]0x201C
$quote2 = [char]0x201D
$sentence2 = "The " + $quote1 + "manage" + $quote2 + " functionality contains the ability to delete a favorite and to contact the landlord to the favorite."
$found2 = $find2.Find.Execute($sentence2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    # The old _GoBack bookmark (if still here) sits inside this range; a
    # plain two-step Text replace collapses the whole range (runs +
    # bookmark) down to one run, which is exactly what the diff wants.
    $find2.Text = "PLACEHOLDERTEXT"
    $find2.Text = $sentence2
}
